$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B19").Value = 30
$ws.Range("B20").Value = 20
$ws.Range("C20").Value = 33
$ws.Range("B21").Value = 12
$ws.Range("B22").Value = 6
$ws.Range("B23").Value = 0

$ws.Range("A4").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C21").Select()
